# The "adductName" column (always "[M-H]-" in this dataset) is redundant
# and was removed from three of the four sheets in this workbook.  In the
# "Corrected" and "Normalized" sheets it was column C; in "PoolAfterDF" it
# was column B.  Deleting the column shifts the remaining columns left,
# which is exactly what Excel's column Delete does.

$wb = $excel.ActiveWorkbook

$wsCorrected = $wb.Worksheets.Item("Corrected")
$wsCorrected.Columns.Item(3).Delete()

$wsNormalized = $wb.Worksheets.Item("Normalized")
$wsNormalized.Columns.Item(3).Delete()

$wsPoolAfterDF = $wb.Worksheets.Item("PoolAfterDF")
$wsPoolAfterDF.Columns.Item(2).Delete()
